# B6-PowerPoint.pptx edit: re-style the three summary tables (slides 14-16)
# from the deck's local custom table style to the built-in themed table
# style "Medium Style 2 - Accent 1" ({EFF24B7C-76A3-4670-AEF2-ABAA5FD5F84C}).

$p = $ppt.ActivePresentation

$oldStyleId = "{82B68193-6186-4F62-B64E-E4511BABF635}"
$newStyleId = "{EFF24B7C-76A3-4670-AEF2-ABAA5FD5F84C}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
